$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete trailing rows (sheet shrinks from 19 to 16 rows)
$ws.Rows("17:19").Delete()

# Rewrite the data rows (2-16) with the refreshed TPM-derived values
$data = New-Object "object[,]" 15,20
$data[0,0] = "ECs"
$data[0,1] = "Lgi2"
$data[0,2] = "Adam23"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 1
$data[0,6] = 0.3777555
$data[0,7] = 0.755511
$data[0,8] = 0.04674878051708328
$data[0,9] = 0.03234814298672928
$data[0,10] = 2
$data[0,11] = 1
$data[0,12] = 0.4562695
$data[0,13] = 0.912539
$data[0,14] = 0.01609359429837405
$data[0,15] = 0.01172153108534722
$data[0,16] = 0.17235831310725
$data[0,17] = 0.689433252429
$data[0,18] = 0.0007523559075856716
$data[0,19] = 0.0003791697635722041

$data[1,0] = "ECs"
$data[1,1] = "Lgi2"
$data[1,2] = "Adam23"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 1
$data[1,6] = 0.3777555
$data[1,7] = 0.755511
$data[1,8] = 0.04674878051708328
$data[1,9] = 0.03234814298672928
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 21.108446
$data[1,13] = 63.325338
$data[1,14] = 0.7445397209174328
$data[1,15] = 0.8134117203287967
$data[1,16] = 7.973831572953001
$data[1,17] = 47.842989437718
$data[1,18] = 0.03480632399941951
$data[1,19] = 0.02631235863627737

$data[2,0] = "ECs"
$data[2,1] = "Lgi2"
$data[2,2] = "Adam23"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 2
$data[2,5] = 1
$data[2,6] = 0.3777555
$data[2,7] = 0.755511
$data[2,8] = 0.04674878051708328
$data[2,9] = 0.03234814298672928
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.006762666666666667
$data[2,13] = 0.020288
$data[2,14] = 0.0002385336160064851
$data[2,15] = 0.0002605986403425218
$data[2,16] = 0.002554634528
$data[2,17] = 0.015327807168
$data[2,18] = [double]"1.115115566063339E-05"
$data[2,19] = [double]"8.429882079947135E-06"

$data[3,0] = "ECs"
$data[3,1] = "Lgi2"
$data[3,2] = "Adam23"
$data[3,3] = "MuSCs"
$data[3,4] = 2
$data[3,5] = 1
$data[3,6] = 0.3777555
$data[3,7] = 0.755511
$data[3,8] = 0.04674878051708328
$data[3,9] = 0.03234814298672928
$data[3,10] = 2
$data[3,11] = 1
$data[3,12] = 6.745213
$data[3,13] = 13.490426
$data[3,14] = 0.2379179881147404
$data[3,15] = 0.1732840434365834
$data[3,16] = 2.5480413094215
$data[3,17] = 10.192165237686
$data[3,18] = 0.01112237580744203
$data[3,19] = 0.005605417014405209

$data[4,0] = "ECs"
$data[4,1] = "Lgi2"
$data[4,2] = "Adam23"
$data[4,3] = "Neutrophils"
$data[4,4] = 2
$data[4,5] = 1
$data[4,6] = 0.3777555
$data[4,7] = 0.755511
$data[4,8] = 0.04674878051708328
$data[4,9] = 0.03234814298672928
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.03430933333333333
$data[4,13] = 0.102928
$data[4,14] = 0.00121016305344615
$data[4,15] = 0.00132210650893016
$data[4,16] = 0.012960539368
$data[4,17] = 0.07776323620799999
$data[4,18] = [double]"5.65736469754374E-05"
$data[4,19] = [double]"4.276769039455829E-05"

$data[5,0] = "FAPs"
$data[5,1] = "Lgi2"
$data[5,2] = "Adam23"
$data[5,3] = "ECs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 7.194537333333333
$data[5,7] = 21.583612
$data[5,8] = 0.8903532753804024
$data[5,9] = 0.9241291882528327
$data[5,10] = 2
$data[5,11] = 1
$data[5,12] = 0.4562695
$data[5,13] = 0.912539
$data[5,14] = 0.01609359429837405
$data[5,15] = 0.01172153108534722
$data[5,16] = 3.282647951811333
$data[5,17] = 19.695887710868
$data[5,18] = 0.01432898439620071
$data[5,19] = 0.01083220900698228

$data[6,0] = "FAPs"
$data[6,1] = "Lgi2"
$data[6,2] = "Adam23"
$data[6,3] = "FAPs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 7.194537333333333
$data[6,7] = 21.583612
$data[6,8] = 0.8903532753804024
$data[6,9] = 0.9241291882528327
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 21.108446
$data[6,13] = 63.325338
$data[6,14] = 0.7445397209174328
$data[6,15] = 0.8134117203287967
$data[6,16] = 151.8655027956507
$data[6,17] = 1366.789525160856
$data[6,18] = 0.662903379169647
$data[6,19] = 0.7516975128227911

$data[7,0] = "FAPs"
$data[7,1] = "Lgi2"
$data[7,2] = "Adam23"
$data[7,3] = "Inflammatory-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 7.194537333333333
$data[7,7] = 21.583612
$data[7,8] = 0.8903532753804024
$data[7,9] = 0.9241291882528327
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.006762666666666667
$data[7,13] = 0.020288
$data[7,14] = 0.0002385336160064851
$data[7,15] = 0.0002605986403425218
$data[7,16] = 0.04865425780622222
$data[7,17] = 0.437888320256
$data[7,18] = 0.0002123791862997052
$data[7,19] = 0.0002408268099595266

$data[8,0] = "FAPs"
$data[8,1] = "Lgi2"
$data[8,2] = "Adam23"
$data[8,3] = "MuSCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 7.194537333333333
$data[8,7] = 21.583612
$data[8,8] = 0.8903532753804024
$data[8,9] = 0.9241291882528327
$data[8,10] = 2
$data[8,11] = 1
$data[8,12] = 6.745213
$data[8,13] = 13.490426
$data[8,14] = 0.2379179881147404
$data[8,15] = 0.1732840434365834
$data[8,16] = 48.52868674978533
$data[8,17] = 291.172120498712
$data[8,18] = 0.2118310599898748
$data[8,19] = 0.1601368423982185

$data[9,0] = "FAPs"
$data[9,1] = "Lgi2"
$data[9,2] = "Adam23"
$data[9,3] = "Neutrophils"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 7.194537333333333
$data[9,7] = 21.583612
$data[9,8] = 0.8903532753804024
$data[9,9] = 0.9241291882528327
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 0.03430933333333333
$data[9,13] = 0.102928
$data[9,14] = 0.00121016305344615
$data[9,15] = 0.00132210650893016
$data[9,16] = 0.2468397795484444
$data[9,17] = 2.221558015936
$data[9,18] = 0.001077472638380129
$data[9,19] = 0.001221797214881415

$data[10,0] = "MuSCs"
$data[10,1] = "Lgi2"
$data[10,2] = "Adam23"
$data[10,3] = "ECs"
$data[10,4] = 2
$data[10,5] = 1
$data[10,6] = 0.5082495
$data[10,7] = 1.016499
$data[10,8] = 0.06289794410251424
$data[10,9] = 0.04352266876043807
$data[10,10] = 2
$data[10,11] = 1
$data[10,12] = 0.4562695
$data[10,13] = 0.912539
$data[10,14] = 0.01609359429837405
$data[10,15] = 0.01172153108534722
$data[10,16] = 0.23189874524025
$data[10,17] = 0.927594980961
$data[10,18] = 0.001012253994587673
$data[10,19] = 0.0005101523147927454

$data[11,0] = "MuSCs"
$data[11,1] = "Lgi2"
$data[11,2] = "Adam23"
$data[11,3] = "FAPs"
$data[11,4] = 2
$data[11,5] = 1
$data[11,6] = 0.5082495
$data[11,7] = 1.016499
$data[11,8] = 0.06289794410251424
$data[11,9] = 0.04352266876043807
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 21.108446
$data[11,13] = 63.325338
$data[11,14] = 0.7445397209174328
$data[11,15] = 0.8134117203287967
$data[11,16] = 10.728357125277
$data[11,17] = 64.370142751662
$data[11,18] = 0.04683001774836625
$data[11,19] = 0.03540184886972831

$data[12,0] = "MuSCs"
$data[12,1] = "Lgi2"
$data[12,2] = "Adam23"
$data[12,3] = "Inflammatory-Mac"
$data[12,4] = 2
$data[12,5] = 1
$data[12,6] = 0.5082495
$data[12,7] = 1.016499
$data[12,8] = 0.06289794410251424
$data[12,9] = 0.04352266876043807
$data[12,10] = 1
$data[12,11] = 0.3333333333333333
$data[12,12] = 0.006762666666666667
$data[12,13] = 0.020288
$data[12,14] = 0.0002385336160064851
$data[12,15] = 0.0002605986403425218
$data[12,16] = 0.003437121952
$data[12,17] = 0.020622731712
$data[12,18] = [double]"1.50032740461465E-05"
$data[12,19] = [double]"1.134194830304811E-05"

$data[13,0] = "MuSCs"
$data[13,1] = "Lgi2"
$data[13,2] = "Adam23"
$data[13,3] = "MuSCs"
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 0.5082495
$data[13,7] = 1.016499
$data[13,8] = 0.06289794410251424
$data[13,9] = 0.04352266876043807
$data[13,10] = 2
$data[13,11] = 1
$data[13,12] = 6.745213
$data[13,13] = 13.490426
$data[13,14] = 0.2379179881147404
$data[13,15] = 0.1732840434365834
$data[13,16] = 3.4282511346435
$data[13,17] = 13.713004538574
$data[13,18] = 0.01496455231742359
$data[13,19] = 0.007541784023959784

$data[14,0] = "MuSCs"
$data[14,1] = "Lgi2"
$data[14,2] = "Adam23"
$data[14,3] = "Neutrophils"
$data[14,4] = 2
$data[14,5] = 1
$data[14,6] = 0.5082495
$data[14,7] = 1.016499
$data[14,8] = 0.06289794410251424
$data[14,9] = 0.04352266876043807
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.03430933333333333
$data[14,13] = 0.102928
$data[14,14] = 0.00121016305344615
$data[14,15] = 0.00132210650893016
$data[14,16] = 0.017437701512
$data[14,17] = 0.104626209072
$data[14,18] = [double]"7.611676809058391E-05"
$data[14,19] = [double]"5.75416036541865E-05"

$ws.Range("A2:T16").Value = $data

Write-Output "rows written"
